$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.380719000000001
$ws.Range("H2").Value = 28.142157
$ws.Range("I2").Value = 0.03679977590837273
$ws.Range("J2").Value = 0.03679977590837273
$ws.Range("M2").Value = 0.243056
$ws.Range("N2").Value = 0.729168
$ws.Range("O2").Value = 0.002199620488481675
$ws.Range("P2").Value = 0.002199620488481675
$ws.Range("Q2").Value = 2.280040037264
$ws.Range("R2").Value = 20.520360335376
$ws.Range("S2").Value = 0.00008094554105959099
$ws.Range("T2").Value = 0.00008094554105959099
$ws.Range("G3").Value = 9.380719000000001
$ws.Range("H3").Value = 28.142157
$ws.Range("I3").Value = 0.03679977590837273
$ws.Range("J3").Value = 0.03679977590837273
$ws.Range("M3").Value = 70.95253000000001
$ws.Range("N3").Value = 212.85759
$ws.Range("O3").Value = 0.6421097964979703
$ws.Range("P3").Value = 0.6421097964979703
$ws.Range("Q3").Value = 665.5857462690701
$ws.Range("R3").Value = 5990.271716421631
$ws.Range("S3").Value = 0.02362949661969612
$ws.Range("T3").Value = 0.02362949661969612
$ws.Range("G4").Value = 9.380719000000001
$ws.Range("H4").Value = 28.142157
$ws.Range("I4").Value = 0.03679977590837273
$ws.Range("J4").Value = 0.03679977590837273
$ws.Range("M4").Value = 0.04794200000000001
$ws.Range("N4").Value = 0.143826
$ws.Range("O4").Value = 0.0004338679376719292
$ws.Range("P4").Value = 0.0004338679376719292
$ws.Range("Q4").Value = 0.4497304302980001
$ws.Range("R4").Value = 4.047573872682
$ws.Range("S4").Value = 0.00001596624288015482
$ws.Range("T4").Value = 0.00001596624288015482
$ws.Range("G5").Value = 9.380719000000001
$ws.Range("H5").Value = 28.142157
$ws.Range("I5").Value = 0.03679977590837273
$ws.Range("J5").Value = 0.03679977590837273
$ws.Range("M5").Value = 39.25553366666666
$ws.Range("N5").Value = 117.766601
$ws.Range("O5").Value = 0.3552567150758761
$ws.Range("P5").Value = 0.3552567150758761
$ws.Range("Q5").Value = 368.2451305220397
$ws.Range("R5").Value = 3314.206174698357
$ws.Range("S5").Value = 0.01307336750473686
$ws.Range("T5").Value = 0.01307336750473686
$ws.Range("I6").Value = 0.3547860986448385
$ws.Range("J6").Value = 0.3547860986448385
$ws.Range("M6").Value = 0.243056
$ws.Range("N6").Value = 0.729168
$ws.Range("O6").Value = 0.002199620488481675
$ws.Range("P6").Value = 0.002199620488481675
$ws.Range("Q6").Value = 21.981832486944
$ws.Range("R6").Value = 197.836492382496
$ws.Range("S6").Value = 0.0007803947716076673
$ws.Range("T6").Value = 0.0007803947716076673
$ws.Range("I7").Value = 0.3547860986448385
$ws.Range("J7").Value = 0.3547860986448385
$ws.Range("M7").Value = 70.95253000000001
$ws.Range("N7").Value = 212.85759
$ws.Range("O7").Value = 0.6421097964979703
$ws.Range("P7").Value = 0.6421097964979703
$ws.Range("Q7").Value = 6416.902396916221
$ws.Range("R7").Value = 57752.12157224599
$ws.Range("S7").Value = 0.2278116296011461
$ws.Range("T7").Value = 0.2278116296011461
$ws.Range("I8").Value = 0.3547860986448385
$ws.Range("J8").Value = 0.3547860986448385
$ws.Range("M8").Value = 0.04794200000000001
$ws.Range("N8").Value = 0.143826
$ws.Range("O8").Value = 0.0004338679376719292
$ws.Range("P8").Value = 0.0004338679376719292
$ws.Range("Q8").Value = 4.335844468308
$ws.Range("R8").Value = 39.02260021477201
$ws.Range("S8").Value = 0.0001539303129337057
$ws.Range("T8").Value = 0.0001539303129337057
$ws.Range("I9").Value = 0.3547860986448385
$ws.Range("J9").Value = 0.3547860986448385
$ws.Range("M9").Value = 39.25553366666666
$ws.Range("N9").Value = 117.766601
$ws.Range("O9").Value = 0.3552567150758761
$ws.Range("P9").Value = 0.3552567150758761
$ws.Range("Q9").Value = 3550.245890849258
$ws.Range("R9").Value = 31952.21301764332
$ws.Range("S9").Value = 0.1260401439591511
$ws.Range("T9").Value = 0.1260401439591511
$ws.Range("G10").Value = 100.179423
$ws.Range("H10").Value = 300.538269
$ws.Range("I10").Value = 0.3929954960840508
$ws.Range("J10").Value = 0.3929954960840508
$ws.Range("M10").Value = 0.243056
$ws.Range("N10").Value = 0.729168
$ws.Range("O10").Value = 0.002199620488481675
$ws.Range("P10").Value = 0.002199620488481675
$ws.Range("Q10").Value = 24.349209836688
$ws.Range("R10").Value = 219.142888530192
$ws.Range("S10").Value = 0.0008644409450674977
$ws.Range("T10").Value = 0.0008644409450674977
$ws.Range("G11").Value = 100.179423
$ws.Range("H11").Value = 300.538269
$ws.Range("I11").Value = 0.3929954960840508
$ws.Range("J11").Value = 0.3929954960840508
$ws.Range("M11").Value = 70.95253000000001
$ws.Range("N11").Value = 212.85759
$ws.Range("O11").Value = 0.6421097964979703
$ws.Range("P11").Value = 0.6421097964979703
$ws.Range("Q11").Value = 7107.983515790191
$ws.Range("R11").Value = 63971.85164211172
$ws.Range("S11").Value = 0.2523462580151487
$ws.Range("T11").Value = 0.2523462580151487
$ws.Range("G12").Value = 100.179423
$ws.Range("H12").Value = 300.538269
$ws.Range("I12").Value = 0.3929954960840508
$ws.Range("J12").Value = 0.3929954960840508
$ws.Range("M12").Value = 0.04794200000000001
$ws.Range("N12").Value = 0.143826
$ws.Range("O12").Value = 0.0004338679376719292
$ws.Range("P12").Value = 0.0004338679376719292
$ws.Range("Q12").Value = 4.802801897466001
$ws.Range("R12").Value = 43.225217077194
$ws.Range("S12").Value = 0.0001705081454003439
$ws.Range("T12").Value = 0.0001705081454003439
$ws.Range("G13").Value = 100.179423
$ws.Range("H13").Value = 300.538269
$ws.Range("I13").Value = 0.3929954960840508
$ws.Range("J13").Value = 0.3929954960840508
$ws.Range("M13").Value = 39.25553366666666
$ws.Range("N13").Value = 117.766601
$ws.Range("O13").Value = 0.3552567150758761
$ws.Range("P13").Value = 0.3552567150758761
$ws.Range("Q13").Value = 3932.596712283741
$ws.Range("R13").Value = 35393.37041055367
$ws.Range("S13").Value = 0.1396142889784342
$ws.Range("T13").Value = 0.1396142889784342
$ws.Range("G14").Value = 1.427630666666667
$ws.Range("H14").Value = 4.282892
$ws.Range("I14").Value = 0.005600475679236752
$ws.Range("J14").Value = 0.005600475679236752
$ws.Range("M14").Value = 0.243056
$ws.Range("N14").Value = 0.729168
$ws.Range("O14").Value = 0.002199620488481675
$ws.Range("P14").Value = 0.002199620488481675
$ws.Range("Q14").Value = 0.3469941993173334
$ws.Range("R14").Value = 3.122947793856
$ws.Range("S14").Value = 0.00001231892104929248
$ws.Range("T14").Value = 0.00001231892104929248
$ws.Range("G15").Value = 1.427630666666667
$ws.Range("H15").Value = 4.282892
$ws.Range("I15").Value = 0.005600475679236752
$ws.Range("J15").Value = 0.005600475679236752
$ws.Range("M15").Value = 70.95253000000001
$ws.Range("N15").Value = 212.85759
$ws.Range("O15").Value = 0.6421097964979703
$ws.Range("P15").Value = 0.6421097964979703
$ws.Range("Q15").Value = 101.2940077055867
$ws.Range("R15").Value = 911.6460693502802
$ws.Range("S15").Value = 0.003596120298686542
$ws.Range("T15").Value = 0.003596120298686542
$ws.Range("G16").Value = 1.427630666666667
$ws.Range("H16").Value = 4.282892
$ws.Range("I16").Value = 0.005600475679236752
$ws.Range("J16").Value = 0.005600475679236752
$ws.Range("M16").Value = 0.04794200000000001
$ws.Range("N16").Value = 0.143826
$ws.Range("O16").Value = 0.0004338679376719292
$ws.Range("P16").Value = 0.0004338679376719292
$ws.Range("Q16").Value = 0.06844346942133335
$ws.Range("R16").Value = 0.6159912247920001
$ws.Range("S16").Value = 0.000002429866832932247
$ws.Range("T16").Value = 0.000002429866832932247
$ws.Range("G17").Value = 1.427630666666667
$ws.Range("H17").Value = 4.282892
$ws.Range("I17").Value = 0.005600475679236752
$ws.Range("J17").Value = 0.005600475679236752
$ws.Range("M17").Value = 39.25553366666666
$ws.Range("N17").Value = 117.766601
$ws.Range("O17").Value = 0.3552567150758761
$ws.Range("P17").Value = 0.3552567150758761
$ws.Range("Q17").Value = 56.04240369889911
$ws.Range("R17").Value = 504.381633290092
$ws.Range("S17").Value = 0.001989606592667984
$ws.Range("T17").Value = 0.001989606592667985
$ws.Range("G18").Value = 53.48524799999999
$ws.Range("H18").Value = 160.455744
$ws.Range("I18").Value = 0.2098181536835013
$ws.Range("J18").Value = 0.2098181536835013
$ws.Range("M18").Value = 0.243056
$ws.Range("N18").Value = 0.729168
$ws.Range("O18").Value = 0.002199620488481675
$ws.Range("P18").Value = 0.002199620488481675
$ws.Range("Q18").Value = 12.999910437888
$ws.Range("R18").Value = 116.999193940992
$ws.Range("S18").Value = 0.0004615203096976261
$ws.Range("T18").Value = 0.0004615203096976261
$ws.Range("G19").Value = 53.48524799999999
$ws.Range("H19").Value = 160.455744
$ws.Range("I19").Value = 0.2098181536835013
$ws.Range("J19").Value = 0.2098181536835013
$ws.Range("M19").Value = 70.95253000000001
$ws.Range("N19").Value = 212.85759
$ws.Range("O19").Value = 0.6421097964979703
$ws.Range("P19").Value = 0.6421097964979703
$ws.Range("Q19").Value = 3794.91366327744
$ws.Range("R19").Value = 34154.22296949696
$ws.Range("S19").Value = 0.1347262919632929
$ws.Range("T19").Value = 0.1347262919632929
$ws.Range("G20").Value = 53.48524799999999
$ws.Range("H20").Value = 160.455744
$ws.Range("I20").Value = 0.2098181536835013
$ws.Range("J20").Value = 0.2098181536835013
$ws.Range("M20").Value = 0.04794200000000001
$ws.Range("N20").Value = 0.143826
$ws.Range("O20").Value = 0.0004338679376719292
$ws.Range("P20").Value = 0.0004338679376719292
$ws.Range("Q20").Value = 2.564189759616
$ws.Range("R20").Value = 23.077707836544
$ws.Range("S20").Value = 0.0000910333696247926
$ws.Range("T20").Value = 0.0000910333696247926
$ws.Range("G21").Value = 53.48524799999999
$ws.Range("H21").Value = 160.455744
$ws.Range("I21").Value = 0.2098181536835013
$ws.Range("J21").Value = 0.2098181536835013
$ws.Range("M21").Value = 39.25553366666666
$ws.Range("N21").Value = 117.766601
$ws.Range("O21").Value = 0.3552567150758761
$ws.Range("P21").Value = 0.3552567150758761
$ws.Range("Q21").Value = 2099.591953534015
$ws.Range("R21").Value = 18896.32758180614
$ws.Range("S21").Value = 0.07453930804088599
$ws.Range("T21").Value = 0.074539308040886
